# Update the division problems in the practice table.
# Each cell is addressed directly by (row, column) and its Range.Text is
# set explicitly. This keeps each replacement scoped strictly to its own
# cell, which matters because some new values coincide with old values
# used elsewhere in the table (e.g. "71÷4=" is both a source value in one
# cell and the target value written into a different cell).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "77÷7="
$t.Cell(1, 2).Range.Text  = "10÷7="
$t.Cell(1, 3).Range.Text  = "35÷6="
$t.Cell(1, 4).Range.Text  = "21÷3="
$t.Cell(1, 5).Range.Text  = "96÷7="

$t.Cell(5, 1).Range.Text  = "86÷4="
$t.Cell(5, 2).Range.Text  = "80÷4="
$t.Cell(5, 3).Range.Text  = "93÷7="
$t.Cell(5, 4).Range.Text  = "56÷3="
$t.Cell(5, 5).Range.Text  = "21÷5="

$t.Cell(9, 1).Range.Text  = "28÷8="
$t.Cell(9, 2).Range.Text  = "98÷9="
$t.Cell(9, 3).Range.Text  = "27÷7="
$t.Cell(9, 4).Range.Text  = "20÷8="
$t.Cell(9, 5).Range.Text  = "50÷6="

$t.Cell(13, 1).Range.Text = "30÷8="
$t.Cell(13, 2).Range.Text = "71÷4="
$t.Cell(13, 3).Range.Text = "73÷7="
$t.Cell(13, 4).Range.Text = "12÷5="
$t.Cell(13, 5).Range.Text = "37÷5="

$t.Cell(17, 1).Range.Text = "48÷5="
$t.Cell(17, 2).Range.Text = "49÷4="
$t.Cell(17, 3).Range.Text = "50÷4="
$t.Cell(17, 4).Range.Text = "44÷6="
$t.Cell(17, 5).Range.Text = "87÷8="

Write-Host "Done updating division problems."
